# Update countries & provincias Spain
# - Refresh case numbers for a set of countries (same underlying data, rows
#   are keyed by the country name in column A, not by row position).
# - After refreshing, the sheet is re-sorted descending by "Casos totales"
#   (column B), which naturally reshuffles a few rows (Indonesia/China,
#   Oman/Filipinas, Eslovenia/Cabo Verde/Guinea-Bisau, ...).
# - The "last updated" timestamp banner in A1 is bumped.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$dataRange = $ws.Range("A4:A219")

# Country => new [Casos totales, Nuevos casos, Casos activos, Recuperados,
#                 Casos criticos, Muertes hoy, Muertes]
$updates = @{
    "Banglades"                          = @(202066, 2709, 110098, 89387, 0, 34, 2581)
    "China"                              = @(83644, 22, 78758, 252, 0, 0, 4634)
    "Indonesia"                          = @(84882, 1752, 43268, 37598, 0, 59, 4016)
    "Filipinas"                          = @(65304, 2357, 22067, 41464, 0, 113, 1773)
    "Oman"                               = @(65504, 1311, 42772, 22424, 0, 10, 308)
    "Polonia"                            = @(39746, 339, 29924, 8204, 0, 6, 1618)
    "Austria"                            = @(19573, 134, 17501, 1361, 0, 0, 711)
    "Malasia"                            = @(8764, 9, 8546, 96, 0, 0, 122)
    "Consejo Danes para los Refugiados"  = @(8324, 75, 4313, 3818, 0, 0, 193)
    "Estado de Palestina"                = @(7764, 0, 1492, 6217, 0, 2, 55)
    "Finlandia"                          = @(7318, 17, 6880, 110, 0, 0, 328)
    "Cabo Verde"                         = @(1939, 0, 902, 1018, 0, 0, 19)
    "Guinea-Bisau"                       = @(1927, 0, 773, 1128, 0, 0, 26)
    "Eslovenia"                          = @(1940, 24, 1568, 261, 0, 0, 111)
    "Macao"                              = @(46, 0, 46, 0, 0, 0, 0)
}

foreach ($country in $updates.Keys) {
    $found = $dataRange.Find($country)
    if ($found -eq $null) {
        continue
    }
    $r = $found.Row
    $vals = $updates[$country]
    for ($i = 0; $i -lt $vals.Length; $i++) {
        $ws.Cells.Item($r, 2 + $i).Value = $vals[$i]
    }
}

# Re-sort the country table (rows 4-219) descending by "Casos totales" (col B)
$sortRange = $ws.Range("A4:H219")
$keyRange = $ws.Range("B4:B219")
$sortRange.Sort($keyRange, 2)

# Bump the "last updated" banner
$ws.Range("A1").Value = "Datos actualizados a 18 de Julio de 2020 a las 11:27"
